$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Rule "R10" (row 10): update the "From" threshold value (column C)
# from 18 to 100, per the commit's saved change.
$ws.Range("C10").Value = 100
